$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H9").Value = 272.47058
$ws.Range("I9").Value = 99.90909000000001
$ws.Range("J9").Value = 588.8333
$ws.Range("K9").Value = 99.90909000000001
$ws.Range("L9").Value = 588.8333
$ws.Range("M9").Value = 69.09090999999999
$ws.Range("N9").Value = -926.8333

$ws.Range("H15").Value = 1279.3
$ws.Range("I15").Value = 1279.3
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 3837.9
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -3668.9

$ws.Range("H116").Value = 1181.7858
$ws.Range("I116").Value = 822.2727
$ws.Range("J116").Value = 2500
$ws.Range("K116").Value = 822.2727
$ws.Range("L116").Value = 2500
$ws.Range("M116").Value = 2619.7273
$ws.Range("N116").Value = -9384

$ws = $wb.Worksheets("ARM")
$ws.Range("H33").Value = 17157.2
$ws.Range("I33").Value = 3542.3333
$ws.Range("J33").Value = 37579.5
$ws.Range("K33").Value = 3542.3333
$ws.Range("L33").Value = 37579.5
$ws.Range("M33").Value = -3213.3333
$ws.Range("N33").Value = -38237.5

$ws.Range("H61").Value = 1409.0735
$ws.Range("I61").Value = 1441.5526
$ws.Range("J61").Value = 1367.9333
$ws.Range("K61").Value = 1441.5526
$ws.Range("L61").Value = 1367.9333
$ws.Range("M61").Value = -1229.5526
$ws.Range("N61").Value = -1791.9333

$ws.Range("H74").Value = 8621709
$ws.Range("I74").Value = 10001051
$ws.Range("J74").Value = 817.75
$ws.Range("K74").Value = 10001051
$ws.Range("L74").Value = 817.75
$ws.Range("M74").Value = -10000177
$ws.Range("N74").Value = -2565.75

$ws.Range("H77").Value = 8621709
$ws.Range("I77").Value = 10001051
$ws.Range("J77").Value = 817.75
$ws.Range("K77").Value = 50005255
$ws.Range("L77").Value = 4088.75
$ws.Range("M77").Value = -50000887
$ws.Range("N77").Value = -12824.75

$ws.Range("H135").Value = 15664.625
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 15664.625
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 15664.625
$ws.Range("N135").Value = -25804.625

$ws.Range("H136").Value = 1409.0735
$ws.Range("I136").Value = 1441.5526
$ws.Range("J136").Value = 1367.9333
$ws.Range("K136").Value = 4324.6578
$ws.Range("L136").Value = 4103.7999
$ws.Range("M136").Value = -1774.6578
$ws.Range("N136").Value = -9203.7999

$ws = $wb.Worksheets("BSM")
$ws.Range("H75").Value = 12663.5
$ws.Range("I75").Value = 6884.6665
$ws.Range("J75").Value = 30000
$ws.Range("K75").Value = 6884.6665
$ws.Range("L75").Value = 30000
$ws.Range("M75").Value = -5948.6665
$ws.Range("N75").Value = -31872

$ws.Range("H78").Value = 12663.5
$ws.Range("I78").Value = 6884.6665
$ws.Range("J78").Value = 30000
$ws.Range("K78").Value = 20653.9995
$ws.Range("L78").Value = 90000
$ws.Range("M78").Value = -15973.9995
$ws.Range("N78").Value = -99360

$ws.Range("H86").Value = 897092.25
$ws.Range("I86").Value = 3146.2354
$ws.Range("J86").Value = 2585657
$ws.Range("K86").Value = 3146.2354
$ws.Range("L86").Value = 2585657
$ws.Range("M86").Value = -2023.2354
$ws.Range("N86").Value = -2587903

$ws.Range("H89").Value = 897092.25
$ws.Range("I89").Value = 3146.2354
$ws.Range("J89").Value = 2585657
$ws.Range("K89").Value = 15731.177
$ws.Range("L89").Value = 12928285
$ws.Range("M89").Value = -10115.177
$ws.Range("N89").Value = -12939517

$ws = $wb.Worksheets("CRP")
$ws.Range("H16").Value = 1642.8387
$ws.Range("I16").Value = 847
$ws.Range("J16").Value = 2609.2144
$ws.Range("K16").Value = 847
$ws.Range("L16").Value = 2609.2144
$ws.Range("M16").Value = -560
$ws.Range("N16").Value = -3183.2144

$ws.Range("H32").Value = 10750
$ws.Range("I32").Value = 1500
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 1500
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -1184
$ws.Range("N32").Value = -20632

$ws.Range("H74").Value = 27614.5
$ws.Range("I74").Value = 5285
$ws.Range("J74").Value = 49944
$ws.Range("K74").Value = 5285
$ws.Range("L74").Value = 49944
$ws.Range("M74").Value = -4411
$ws.Range("N74").Value = -51692

$ws.Range("H77").Value = 27614.5
$ws.Range("I77").Value = 5285
$ws.Range("J77").Value = 49944
$ws.Range("K77").Value = 15855
$ws.Range("L77").Value = 149832
$ws.Range("M77").Value = -11487
$ws.Range("N77").Value = -158568

$ws.Range("H113").Value = 1642.8387
$ws.Range("I113").Value = 847
$ws.Range("J113").Value = 2609.2144
$ws.Range("K113").Value = 847
$ws.Range("L113").Value = 2609.2144
$ws.Range("M113").Value = 1323
$ws.Range("N113").Value = -6949.2144

$ws = $wb.Worksheets("CUL")
$ws.Range("H39").Value = 1077.2727
$ws.Range("I39").Value = 375
$ws.Range("J39").Value = 1478.5714
$ws.Range("K39").Value = 1125
$ws.Range("L39").Value = 4435.7142
$ws.Range("M39").Value = -831
$ws.Range("N39").Value = -5023.7142

$ws.Range("H62").Value = 1745.8334
$ws.Range("I62").Value = 1000
$ws.Range("J62").Value = 2491.6667
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 7475.000100000001
$ws.Range("M62").Value = -2314
$ws.Range("N62").Value = -8847.000100000001

$ws.Range("H65").Value = 1745.8334
$ws.Range("I65").Value = 1000
$ws.Range("J65").Value = 2491.6667
$ws.Range("K65").Value = 9000
$ws.Range("L65").Value = 22425.0003
$ws.Range("M65").Value = -5568
$ws.Range("N65").Value = -29289.0003

$ws = $wb.Worksheets("GSM")
$ws.Range("H103").Value = 14000
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 14000
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 14000
$ws.Range("N103").Value = -16344

$ws = $wb.Worksheets("LTW")
$ws.Range("H61").Value = 1675.909
$ws.Range("I61").Value = 1609.0294
$ws.Range("J61").Value = 1903.3
$ws.Range("K61").Value = 1609.0294
$ws.Range("L61").Value = 1903.3
$ws.Range("M61").Value = -1407.0294
$ws.Range("N61").Value = -2307.3

$ws.Range("H113").Value = 1675.909
$ws.Range("I113").Value = 1609.0294
$ws.Range("J113").Value = 1903.3
$ws.Range("K113").Value = 1609.0294
$ws.Range("L113").Value = 1903.3
$ws.Range("M113").Value = 560.9706000000001
$ws.Range("N113").Value = -6243.3

$ws.Range("H134").Value = 10000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 10000
$ws.Range("N134").Value = -20140

$ws.Range("H138").Value = 33440
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 33440
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 33440
$ws.Range("N138").Value = -43720

$ws.Range("H141").Value = 10000
$ws.Range("I141").Value = 10000
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 10000
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -4820
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets("WVR")
$ws.Range("H107").Value = 966.25
$ws.Range("I107").Value = 2325.8
$ws.Range("J107").Value = 513.06665
$ws.Range("K107").Value = 6977.400000000001
$ws.Range("L107").Value = 1539.19995
$ws.Range("M107").Value = -5057.400000000001
$ws.Range("N107").Value = -5379.19995

$ws.Range("H133").Value = 25000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 25000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 25000
$ws.Range("N133").Value = -35120

$ws.Range("H136").Value = 1996.0444
$ws.Range("I136").Value = 2248.5833
$ws.Range("J136").Value = 1490.9667
$ws.Range("K136").Value = 6745.749899999999
$ws.Range("L136").Value = 4472.9001
$ws.Range("M136").Value = -4195.749899999999
$ws.Range("N136").Value = -9572.900099999999

Write-Host "Gungnir profit data refreshed"
